$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 7.810333919920698
$ws.Range("D2").Value = 4.670121255836118
$ws.Range("E2").Value = 12.62814539877288
$ws.Range("F2").Value = 23.6761467110982
$ws.Range("G2").Value = 3.623872999658669
$ws.Range("K2").Value = 10.61131047499369
$ws.Range("M2").Value = 13.50365642178098
$ws.Range("O2").Value = 21.13064648655152

# Row 3
$ws.Range("B3").Value = 7.73932860670066
$ws.Range("D3").Value = 4.627492141058521
$ws.Range("E3").Value = 12.50573201948791
$ws.Range("F3").Value = 23.67862014319258
$ws.Range("G3").Value = 3.625934337218601
$ws.Range("K3").Value = 10.19064380165904
$ws.Range("M3").Value = 13.24624243000257
$ws.Range("O3").Value = 21.19105739484059

# Row 4
$ws.Range("B4").Value = 7.697271464954541
$ws.Range("D4").Value = 4.600815314209901
$ws.Range("E4").Value = 12.43493596842788
$ws.Range("F4").Value = 23.68788812236721
$ws.Range("G4").Value = 3.627266792985701
$ws.Range("K4").Value = 9.921425676798517
$ws.Range("M4").Value = 13.08857049185581
$ws.Range("O4").Value = 21.23401836597061

# Row 5
$ws.Range("B5").Value = 7.680539020038377
$ws.Range("D5").Value = 4.589823146800896
$ws.Range("E5").Value = 12.40721548350868
$ws.Range("F5").Value = 23.69361024634976
$ws.Range("G5").Value = 3.627826627991773
$ws.Range("K5").Value = 9.809069939572721
$ws.Range("M5").Value = 13.02450266791046
$ws.Range("O5").Value = 21.25299516830934

# Row 6
$ws.Range("B6").Value = 7.677785678071267
$ws.Range("D6").Value = 4.587990741625355
$ws.Range("E6").Value = 12.40268160212567
$ws.Range("F6").Value = 23.69467776857839
$ws.Range("G6").Value = 3.627920607377865
$ws.Range("K6").Value = 9.790256676344036
$ws.Range("M6").Value = 13.0138780348347
$ws.Range("O6").Value = 21.25623485099378

# Row 7
$ws.Range("B7").Value = 7.697044136413306
$ws.Range("D7").Value = 4.600667553623904
$ws.Range("E7").Value = 12.43455750747695
$ws.Range("F7").Value = 23.68795742182191
$ws.Range("G7").Value = 3.627274274825552
$ws.Range("K7").Value = 9.919920983572817
$ws.Range("M7").Value = 13.08770558433181
$ws.Range("O7").Value = 21.23426835034035

# Row 8
$ws.Range("B8").Value = 7.785543101281153
$ws.Range("D8").Value = 4.655530091863184
$ws.Range("E8").Value = 12.58505083084777
$ws.Range("F8").Value = 23.67538968710507
$ws.Range("G8").Value = 3.624569920260346
$ws.Range("K8").Value = 10.46858875520209
$ws.Range("M8").Value = 13.414873848409
$ws.Range("O8").Value = 21.15025478937456

# Row 9
$ws.Range("B9").Value = 7.97045032340441
$ws.Range("D9").Value = 4.75890646769477
$ws.Range("E9").Value = 12.91329313320009
$ws.Range("F9").Value = 23.71230700142138
$ws.Range("G9").Value = 3.619794116976586
$ws.Range("K9").Value = 11.45405680522614
$ws.Range("M9").Value = 14.05565218924631
$ws.Range("O9").Value = 21.03231220238075

# Row 10
$ws.Range("B10").Value = 8.112035000801859
$ws.Range("D10").Value = 4.832000571096542
$ws.Range("E10").Value = 13.17244973476686
$ws.Range("F10").Value = 23.77695987451093
$ws.Range("G10").Value = 3.616603379812751
$ws.Range("K10").Value = 12.11870008007091
$ws.Range("M10").Value = 14.52099331408035
$ws.Range("O10").Value = 20.97450598011496

# Row 11
$ws.Range("B11").Value = 8.177429290007133
$ws.Range("D11").Value = 4.864572995945374
$ws.Range("E11").Value = 13.29375008897548
$ws.Range("F11").Value = 23.81448606556001
$ws.Range("G11").Value = 3.615220148726813
$ws.Range("K11").Value = 12.40747928166741
$ws.Range("M11").Value = 14.73052144936299
$ws.Range("O11").Value = 20.9545299680001

# Row 12
$ws.Range("B12").Value = 8.202312172244183
$ws.Range("D12").Value = 4.876804853292209
$ws.Range("E12").Value = 13.34012955915126
$ws.Range("F12").Value = 23.8298574024854
$ws.Range("G12").Value = 3.614706113193372
$ws.Range("K12").Value = 12.51483463346793
$ws.Range("M12").Value = 14.80947687001712
$ws.Range("O12").Value = 20.94787840969867

# Row 13
$ws.Range("B13").Value = 8.196948240561078
$ws.Range("D13").Value = 4.874175146759281
$ws.Range("E13").Value = 13.33012180535168
$ws.Range("F13").Value = 23.82649538871948
$ws.Range("G13").Value = 3.614816386482215
$ws.Range("K13").Value = 12.49180332692367
$ws.Range("M13").Value = 14.79249087601116
$ws.Range("O13").Value = 20.94927027978848

# Row 14
$ws.Range("B14").Value = 8.179474153265261
$ws.Range("D14").Value = 4.865581403208333
$ws.Range("E14").Value = 13.29755708157592
$ws.Range("F14").Value = 23.8157274473085
$ws.Range("G14").Value = 3.615177663315462
$ws.Range("K14").Value = 12.41635178194752
$ws.Range("M14").Value = 14.73702537402364
$ws.Range("O14").Value = 20.95396442241485

# Row 15
$ws.Range("B15").Value = 8.168785666197165
$ws.Range("D15").Value = 4.860303976154369
$ws.Range("E15").Value = 13.27766695312312
$ws.Range("F15").Value = 23.80928275933345
$ws.Range("G15").Value = 3.615400225886606
$ws.Range("K15").Value = 12.36987390744017
$ws.Range("M15").Value = 14.7029983361613
$ws.Range("O15").Value = 20.95695872218528

# Row 16
$ws.Range("B16").Value = 8.107779269974298
$ws.Range("D16").Value = 4.829857799526931
$ws.Range("E16").Value = 13.16458727547292
$ws.Range("F16").Value = 23.77467027786859
$ws.Range("G16").Value = 3.616695146291188
$ws.Range("K16").Value = 12.09955036265981
$ws.Range("M16").Value = 14.50724996854831
$ws.Range("O16").Value = 20.97593899976133

# Row 17
$ws.Range("B17").Value = 8.070590603235061
$ws.Range("D17").Value = 4.81100265362749
$ws.Range("E17").Value = 13.09605819019162
$ws.Range("F17").Value = 23.75551145240906
$ws.Range("G17").Value = 3.617506982831623
$ws.Range("K17").Value = 11.93020358438021
$ws.Range("M17").Value = 14.38655514794496
$ws.Range("O17").Value = 20.98920462459451

# Row 18
$ws.Range("B18").Value = 8.049294883859105
$ws.Range("D18").Value = 4.800094107878205
$ws.Range("E18").Value = 13.0569656419063
$ws.Range("F18").Value = 23.74525611959747
$ws.Range("G18").Value = 3.617980356856376
$ws.Range("K18").Value = 11.83152482484608
$ws.Range("M18").Value = 14.31693597836464
$ws.Range("O18").Value = 20.99742934774484

# Row 19
$ws.Range("B19").Value = 8.042101383908445
$ws.Range("D19").Value = 4.796389900517979
$ws.Range("E19").Value = 13.04378649374059
$ws.Range("F19").Value = 23.7419152661943
$ws.Range("G19").Value = 3.618141738572596
$ws.Range("K19").Value = 11.79789649179477
$ws.Range("M19").Value = 14.29333244334867
$ws.Range("O19").Value = 21.00031610365776

# Row 20
$ws.Range("B20").Value = 8.07453981141213
$ws.Range("D20").Value = 4.813016431786887
$ws.Range("E20").Value = 13.10332006025461
$ws.Range("F20").Value = 23.75747187595892
$ws.Range("G20").Value = 3.617419896629497
$ws.Range("K20").Value = 11.94836316190331
$ws.Range("M20").Value = 14.39942446115658
$ws.Range("O20").Value = 20.98773089726815

# Row 21
$ws.Range("B21").Value = 8.184603651890002
$ws.Range("D21").Value = 4.868108418514015
$ws.Range("E21").Value = 13.30711039334052
$ws.Range("F21").Value = 23.81885880025322
$ws.Range("G21").Value = 3.615071282957796
$ws.Range("K21").Value = 12.4385683267906
$ws.Range("M21").Value = 14.75332807659477
$ws.Range("O21").Value = 20.9525608324966

# Row 22
$ws.Range("B22").Value = 8.257222838149779
$ws.Range("D22").Value = 4.903513859286674
$ws.Range("E22").Value = 13.44287575569237
$ws.Range("F22").Value = 23.86574180633866
$ws.Range("G22").Value = 3.613593217441117
$ws.Range("K22").Value = 12.74727660743911
$ws.Range("M22").Value = 14.98232787498014
$ws.Range("O22").Value = 20.93489791215674

# Row 23
$ws.Range("B23").Value = 8.218409219722908
$ws.Range("D23").Value = 4.88467387715712
$ws.Range("E23").Value = 13.37019432830507
$ws.Range("F23").Value = 23.84010306273787
$ws.Range("G23").Value = 3.614376900004724
$ws.Range("K23").Value = 12.58359478313889
$ws.Range("M23").Value = 14.86034118900674
$ws.Range("O23").Value = 20.94383669913239

# Row 24
$ws.Range("B24").Value = 8.072754109363292
$ws.Range("D24").Value = 4.812106215556767
$ws.Range("E24").Value = 13.10003601309691
$ws.Range("F24").Value = 23.75658320273073
$ws.Range("G24").Value = 3.617459247618859
$ws.Range("K24").Value = 11.94015732179708
$ws.Range("M24").Value = 14.3936069554078
$ws.Range("O24").Value = 20.98839530649119

# Row 25
$ws.Range("B25").Value = 7.919333207619863
$ws.Range("D25").Value = 4.731421373684638
$ws.Range("E25").Value = 12.82116773037839
$ws.Range("F25").Value = 23.69572255851478
$ws.Range("G25").Value = 3.621029995514723
$ws.Range("K25").Value = 11.19762445359406
$ws.Range("M25").Value = 13.88291278595871
$ws.Range("O25").Value = 21.05917393953174
